$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 218, shifting existing rows 218-249 down to 219-250.
$ws.Rows.Item(218).Insert()

# Populate the newly inserted row 218 with the new price-report record.
$ws.Cells.Item(218, 1).Value = 4
$ws.Cells.Item(218, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(218, 3).Value = "Los Lagos"
$ws.Cells.Item(218, 4).Value = 44491
$ws.Cells.Item(218, 5).Value = 10
$ws.Cells.Item(218, 6).Value = 100114001
$ws.Cells.Item(218, 7).Value = "Papa"
$ws.Cells.Item(218, 8).Value = "Asterix"
$ws.Cells.Item(218, 9).Value = "1a (guarda)"
$ws.Cells.Item(218, 10).Value = 350
$ws.Cells.Item(218, 11).Value = 9000
$ws.Cells.Item(218, 12).Value = 9000
$ws.Cells.Item(218, 13).Value = 9000
$ws.Cells.Item(218, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(218, 15).Value = "Provincia de Llanquihue"
$ws.Cells.Item(218, 16).Value = 360
$ws.Cells.Item(218, 17).Value = 25
$ws.Cells.Item(218, 18).Value = "Hortaliza"
